$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 431 in the underlying data
# source; this pushes every existing record from row 431 onward down by
# one row (old row 431 -> new row 432, ..., old row 543 -> new row 544).
$ws.Rows("431:431").Insert(-4121)

# Populate the newly inserted row 431 with the new record's values.
$ws.Range("A431").Value = 3
$ws.Range("B431").Value = "Femacal de La Calera"
$ws.Range("C431").Value = "Coquimbo"
$ws.Range("D431").Value = 44932
$ws.Range("E431").Value = 5
$ws.Range("F431").Value = "Fruta"
$ws.Range("G431").Value = 100108
$ws.Range("H431").Value = "Tropicales y subtropicales"
$ws.Range("I431").Value = 100108002
$ws.Range("J431").Value = "Mango"
$ws.Range("K431").Value = "Sin especificar"
$ws.Range("L431").Value = "Primera"
$ws.Range("M431").Value = 228
$ws.Range("N431").Value = 7000
$ws.Range("O431").Value = 7000
$ws.Range("P431").Value = 7000
$ws.Range("Q431").Value = "$/bandeja 4 kilos"
$ws.Range("R431").Value = "Perú"
$ws.Range("S431").Value = 1750
$ws.Range("T431").Value = 4
